$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 4 with new data (was previously blank, kept style "3")
$ws.Range("A4").Value = "item in dropdown"
$ws.Range("B4").Value = "option"
$ws.Range("C4").Value = "value"
$ws.Range("D4").Value = "id,name,class etc.."

# Update the active selection to D9 as in the final workbook
$ws.Range("D9").Select()
